# The document's footers/header each contain a single inline picture
# (logo). The edit renames those pictures ("image2.png" <-> "image1.png"
# for the two Pearson logos in the footers, and "image1.jpg" ->
# "image2.jpg" for the BTec logo in the first-page header) while leaving
# everything else (size, description/alt text, embedded image data)
# untouched.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Default (odd/all) footer -> word/footer1.xml : "image2.png" -> "image1.png"
$footer1 = $sec.Footers.Item(1)
if ($footer1.Exists -and $footer1.Range.InlineShapes.Count -ge 1) {
    $pearsonShape1 = $footer1.Range.InlineShapes.Item(1)
    $pearsonShape1.Name = "image1.png"
}

# --- First-page footer -> word/footer2.xml : "image2.png" -> "image1.png"
$footer2 = $sec.Footers.Item(2)
if ($footer2.Exists -and $footer2.Range.InlineShapes.Count -ge 1) {
    $pearsonShape2 = $footer2.Range.InlineShapes.Item(1)
    $pearsonShape2.Name = "image1.png"
}

# --- First-page header -> word/header2.xml : "image1.jpg" -> "image2.jpg"
$header2 = $sec.Headers.Item(2)
if ($header2.Exists -and $header2.Range.InlineShapes.Count -ge 1) {
    $btecShape = $header2.Range.InlineShapes.Item(1)
    $btecShape.Name = "image2.jpg"
}
